# Generate Report for Handoff
# Updates the Overview / zh-cn / de-de sheets to reflect the new handoff
# status: drops the c8d950fd-... row (handed back) from every sheet and
# flips the 5bbc677a-... row's status from "Handed back: in sync with
# en-US" to "Ready for handoff" with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-38-12 04:38:36"
$overview.Rows.Item(3).Delete()

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-12 04:38:33"
$zhcn.Rows.Item(3).Delete()

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-12 04:38:36"
$dede.Rows.Item(3).Delete()
